# Update the "Export" sheet data to reflect the new Saldo upload.
# Rows 2-15 are rewritten in place: a new account (EDSON) is inserted at
# the top, a few accounts (LARA/004452597, RODRIGO-13450/005142624,
# EDUARDO/004946997) are removed, a couple of new accounts (JOSE/002687737,
# PAULO/004404248) are added, and PEDRO's (004458624) balance is updated
# and the row moves up in the (descending-by-Saldo) ordering. Everything
# from row 16 (CRISTINA) onward is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("008103455", "EDSON",    52000),
    @("008035153", "CLAUDIO",  42500),
    @("002687737", "JOSE",     18592.76),
    @("004458624", "PEDRO",    12881.86),
    @("004313254", "GUSTAVO",  4292),
    @("004368468", "AHMAD",    3180.45),
    @("004213139", "LEONARDO", 2610.01),
    @("004404248", "PAULO",    1108.48),
    @("004211368", "ILTON",    986.85),
    @("008008723", "REDRAU",   956.58),
    @("004392159", "RODRIGO",  900.21),
    @("003301389", "EDMUNDO",  832.22),
    @("005599726", "JORGE",    787.13),
    @("005685353", "CARLOS",   767.05)
)

$row = 2
foreach ($entry in $data) {
    $acctCell = $ws.Cells.Item($row, 1)
    $acctCell.NumberFormat = "@"
    $acctCell.Value = $entry[0]

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]

    $row++
}
